$d = $word.ActiveDocument

# 1. Remove the "FINALS WEEK START" run from its original location
#    (the paragraph itself is left in place, now empty).
$rng = $d.Range(0, $d.Content.End)
$found = $rng.Find.Execute("FINALS WEEK START", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = ""
}

# 2. Add "FINALS WEEK START" as a new paragraph in the cell holding day "14"
#    and "Final exam period 10:20-12:10" as a new paragraph in the cell
#    holding day "15" (December calendar table, week row containing 13-19).
$t = $d.Tables.Item(5)
$row = $t.Rows.Item(4)

$cell14 = $row.Cells.Item(2)
$newPara14 = $cell14.Range.Paragraphs.Add()
$newPara14.Range.InsertAfter("FINALS WEEK START")

$cell15 = $row.Cells.Item(3)
$newPara15 = $cell15.Range.Paragraphs.Add()
$newPara15.Range.InsertAfter("Final exam period 10:20-12:10")
